$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (experiment 001 parameters)
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 0.8
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = "'1,10"
$ws.Range("I2").Value = "'1,10"
$ws.Range("J2").Value = 0.3

# Clear out rows 3 and 4 (experiments 002 and 003 removed), keep A3/A4/H3/H4/I3/I4 blank but formatted
$ws.Range("A3:J3").ClearContents()
$ws.Range("A4:J4").ClearContents()

# Update the active selection to B2 as per the new view state
$ws.Range("B2").Select()
